$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column E header cells (row 1) -------------------------------------
# Snapshot an untouched, identically-styled cell (A1, style index 1) so
# we can restore the original look of the header row after the merge
# operation below.
$ws.Range("A1").Copy()

# Strip the outside border first so growing the merge doesn't make Excel
# redistribute it into separate left/middle/right-edge border pieces.
$ws.Range("B1:D1").Borders.LineStyle = 0

# Grow the merged title cell from B1:D1 to B1:E1.
$ws.Range("B1:D1").MergeCells = $false
$ws.Range("B1:E1").MergeCells = $true

# Re-apply the original formatting (border/font/alignment) across the
# whole merged range so every cell in it matches the sheet's existing
# look (same style as A1/B1 previously had).
$ws.Range("B1:E1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 2: new "runtime" header ---------------------------------------
$ws.Range("E2").Value = "runtime"
$ws.Range("D2").Copy()
$ws.Range("E2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Rows 4-23: per-model runtime (mm:ss) column ------------------------
$ws.Range("E4").Value = "49:06"
$ws.Range("E5").Value = "29:45"
$ws.Range("E6").Value = "38:31"
$ws.Range("E7").Value = "68:43"
$ws.Range("E8").Value = "235:38"
$ws.Range("E9").Value = "141:40"
$ws.Range("E10").Value = "153:06"
$ws.Range("E11").Value = "381:42"
$ws.Range("E12").Value = "213:38"
$ws.Range("E13").Value = "123:00"
$ws.Range("E14").Value = "141:50"
$ws.Range("E15").Value = "389:56"
$ws.Range("E16").Value = "214:06"
$ws.Range("E17").Value = "132:17"
$ws.Range("E18").Value = "148:47"
$ws.Range("E19").Value = "362:45"
$ws.Range("E20").Value = "119:40"
$ws.Range("E21").Value = "72:03"
$ws.Range("E22").Value = "77:18"
$ws.Range("E23").Value = "195:47"
